$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column layout change -------------------------------------------------
# Before: A | B(Relevance) | C(Topic) | D(Unit) | E(Scale) | F(Time) | G(Principle) | H(30 word explanation) | I(Notes)
# After:  A | B(Relevance) | C(Topic) | D(Unit) | E(Shape)             | F(Principle) | G(30 word explanation) | H(Notes)
#
# The "Time" column (F) is removed outright, and the old "Scale" column (E)
# is repurposed into the new "Shape" column. This naturally shifts the old
# G/H/I (Principle / 30 word explanation / Notes) down into F/G/H.
$ws.Columns.Item(6).Delete()

# --- Relevance column (B) becomes numeric ---------------------------------
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 0
$ws.Range("B11").Value = 0
$ws.Range("B12").Value = 0
$ws.Range("B13").Value = 0
$ws.Range("B14").Value = 0
$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 0

# --- Row 5 annotation content re-done --------------------------------------
$ws.Range("C5").Value = "urgency"

# Relabel the repurposed column E header (after "urgency" above, to match
# shared-string allocation order of the original authored workbook).
$ws.Range("E1").Value = "Shape"

$ws.Range("D5").Value = "n.a."
$ws.Range("E5").Value = "n.a."
$ws.Range("F5").Value = "utilitarian"
$ws.Range("G5").Value = "Calling for increased urgency, motivated by utilitarian idea in the benefit of all. No specific distribution highlighted. "
$ws.Range("H5").ClearContents()

# --- View state -------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B16").Select()
